$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 218-223 with revised monthly export figures
# Row 218
$ws.Range("B218").Value = 7177
$ws.Range("C218").Value = 2135
$ws.Range("D218").Value = 484
$ws.Range("E218").Value = 1478
$ws.Range("N218").Value = 2502
$ws.Range("O218").Value = 2110
$ws.Range("T218").Value = 57
$ws.Range("U218").Value = 48
$ws.Range("X218").Value = 282
$ws.Range("Z218").Value = 156
$ws.Range("AC218").Value = 1759
$ws.Range("AD218").Value = 1500

# Row 219
$ws.Range("B219").Value = 7294
$ws.Range("C219").Value = 1575
$ws.Range("D219").Value = 453
$ws.Range("E219").Value = 935
$ws.Range("N219").Value = 2978
$ws.Range("O219").Value = 2675
$ws.Range("Q219").Value = 91
$ws.Range("T219").Value = 51
$ws.Range("U219").Value = 37
$ws.Range("Z219").Value = 141
$ws.Range("AC219").Value = 2068
$ws.Range("AD219").Value = 1865
$ws.Range("AF219").Value = 32
$ws.Range("AI219").Value = 75

# Row 220
$ws.Range("B220").Value = 7761
$ws.Range("C220").Value = 1496
$ws.Range("D220").Value = 550
$ws.Range("E220").Value = 706
$ws.Range("K220").Value = 157
$ws.Range("N220").Value = 3209
$ws.Range("O220").Value = 2909
$ws.Range("S220").Value = 36
$ws.Range("T220").Value = 49
$ws.Range("U220").Value = 38
$ws.Range("Z220").Value = 147
$ws.Range("AC220").Value = 2297
$ws.Range("AD220").Value = 2075

# Row 221
$ws.Range("B221").Value = 8160
$ws.Range("C221").Value = 1394
$ws.Range("D221").Value = 502
$ws.Range("N221").Value = 3519
$ws.Range("O221").Value = 3004
$ws.Range("T221").Value = 48
$ws.Range("U221").Value = 48
$ws.Range("Z221").Value = 277
$ws.Range("AC221").Value = 2293
$ws.Range("AD221").Value = 2048
$ws.Range("AH221").Value = 14
$ws.Range("AJ221").Value = 145
$ws.Range("AK221").Value = 56
$ws.Range("AL221").Value = 89
$ws.Range("AM221").Value = 46
$ws.Range("AN221").Value = 14
$ws.Range("AO221").Value = 32

# Row 222
$ws.Range("B222").Value = 7621
$ws.Range("N222").Value = 3331
$ws.Range("O222").Value = 2964
$ws.Range("T222").Value = 52
$ws.Range("U222").Value = 38
$ws.Range("AC222").Value = 2312
$ws.Range("AD222").Value = 2073
$ws.Range("AG222").Value = 34
$ws.Range("AH222").Value = 15
$ws.Range("AJ222").Value = 93
$ws.Range("AK222").Value = 36
$ws.Range("AL222").Value = 57
$ws.Range("AM222").Value = 55
$ws.Range("AO222").Value = 41

# Row 223
$ws.Range("B223").Value = 7453
$ws.Range("C223").Value = 1209
$ws.Range("G223").Value = 129
$ws.Range("N223").Value = 3344
$ws.Range("O223").Value = 2947
$ws.Range("Q223").Value = 132
$ws.Range("T223").Value = 37
$ws.Range("U223").Value = 26
$ws.Range("X223").Value = 311
$ws.Range("Z223").Value = 184
$ws.Range("AC223").Value = 2098
$ws.Range("AD223").Value = 1850
$ws.Range("AH223").Value = 15
$ws.Range("AJ223").Value = 105
$ws.Range("AK223").Value = 43
$ws.Range("AM223").Value = 48
$ws.Range("AN223").Value = 15
$ws.Range("AO223").Value = 32

# Add new row 224 for period 01-07-2021
$ws.Range("A224").NumberFormat = "@"
$ws.Range("A224").Value = "01-07-2021"
$ws.Range("A224").Style = "Normal"
$ws.Range("B224").Value = 7943
$ws.Range("C224").Value = 1231
$ws.Range("D224").Value = 461
$ws.Range("E224").Value = 499
$ws.Range("F224").Value = 59
$ws.Range("G224").Value = 134
$ws.Range("H224").Value = 30
$ws.Range("I224").Value = 22
$ws.Range("J224").Value = 24
$ws.Range("K224").Value = 167
$ws.Range("L224").Value = 165
$ws.Range("M224").Value = 2
$ws.Range("N224").Value = 3599
$ws.Range("O224").Value = 3079
$ws.Range("P224").Value = 303
$ws.Range("Q224").Value = 135
$ws.Range("R224").Value = 64
$ws.Range("S224").Value = 17
$ws.Range("T224").Value = 49
$ws.Range("U224").Value = 37
$ws.Range("V224").Value = 11
$ws.Range("W224").Value = 34
$ws.Range("X224").Value = 282
$ws.Range("Y224").Value = 9
$ws.Range("Z224").Value = 168
$ws.Range("AA224").Value = 33
$ws.Range("AB224").Value = 71
$ws.Range("AC224").Value = 2335
$ws.Range("AD224").Value = 2068
$ws.Range("AE224").Value = 98
$ws.Range("AF224").Value = 43
$ws.Range("AG224").Value = 38
$ws.Range("AH224").Value = 13
$ws.Range("AI224").Value = 75
$ws.Range("AJ224").Value = 127
$ws.Range("AK224").Value = 58
$ws.Range("AL224").Value = 70
$ws.Range("AM224").Value = 50
$ws.Range("AN224").Value = 17
$ws.Range("AO224").Value = 33
$ws.Range("AP224").Value = 70
